$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 65.17856999999999
$ws.Range("I33").Value = 64.04000000000001
$ws.Range("J33").Value = 74.666664
$ws.Range("K33").Value = 64.04000000000001
$ws.Range("L33").Value = 74.666664
$ws.Range("M33").Value = 164.96
$ws.Range("N33").Value = -532.666664
# Row 40
$ws.Range("H40").Value = 1223.579
$ws.Range("I40").Value = 1220.8462
$ws.Range("J40").Value = 1229.5
$ws.Range("K40").Value = 1220.8462
$ws.Range("L40").Value = 1229.5
$ws.Range("M40").Value = -1045.8462
$ws.Range("N40").Value = -1579.5
# Row 41
$ws.Range("H41").Value = 537.9167
$ws.Range("I41").Value = 599.8570999999999
$ws.Range("J41").Value = 451.2
$ws.Range("K41").Value = 599.8570999999999
$ws.Range("L41").Value = 451.2
$ws.Range("M41").Value = -159.8570999999999
$ws.Range("N41").Value = -1331.2
# Row 43
$ws.Range("H43").Value = 6267.316
$ws.Range("J43").Value = 984.375
$ws.Range("L43").Value = 984.375
$ws.Range("N43").Value = -1122.375
# Row 111
$ws.Range("H111").Value = 2122.25
$ws.Range("I111").Value = 1946.7
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 5840.1
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -2773.1
$ws.Range("N111").Value = -15134
# Row 137
$ws.Range("H137").Value = 2578092.2
$ws.Range("I137").Value = 879
$ws.Range("J137").Value = 5584841
$ws.Range("K137").Value = 2637
$ws.Range("L137").Value = 16754523
$ws.Range("M137").Value = -87
$ws.Range("N137").Value = -16759623

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 648.43475
$ws.Range("I2").Value = 686.0625
$ws.Range("J2").Value = 562.4286
$ws.Range("K2").Value = 686.0625
$ws.Range("L2").Value = 562.4286
$ws.Range("M2").Value = -573.0625
$ws.Range("N2").Value = -788.4286
# Row 61
$ws.Range("H61").Value = 4698.915
$ws.Range("I61").Value = 4931.6665
$ws.Range("J61").Value = 3564.25
$ws.Range("K61").Value = 4931.6665
$ws.Range("L61").Value = 3564.25
$ws.Range("M61").Value = -4719.6665
$ws.Range("N61").Value = -3988.25
# Row 74
$ws.Range("H74").Value = 3058.36
$ws.Range("I74").Value = 883.25714
$ws.Range("J74").Value = 8133.6
$ws.Range("K74").Value = 883.25714
$ws.Range("L74").Value = 8133.6
$ws.Range("M74").Value = -9.257140000000049
$ws.Range("N74").Value = -9881.6
# Row 77
$ws.Range("H77").Value = 3058.36
$ws.Range("I77").Value = 883.25714
$ws.Range("J77").Value = 8133.6
$ws.Range("K77").Value = 4416.2857
$ws.Range("L77").Value = 40668
$ws.Range("M77").Value = -48.28570000000036
$ws.Range("N77").Value = -49404
# Row 116
$ws.Range("H116").Value = 648.43475
$ws.Range("I116").Value = 686.0625
$ws.Range("J116").Value = 562.4286
$ws.Range("K116").Value = 686.0625
$ws.Range("L116").Value = 562.4286
$ws.Range("M116").Value = 1607.9375
$ws.Range("N116").Value = -5150.4286
# Row 136
$ws.Range("H136").Value = 4698.915
$ws.Range("I136").Value = 4931.6665
$ws.Range("J136").Value = 3564.25
$ws.Range("K136").Value = 14794.9995
$ws.Range("L136").Value = 10692.75
$ws.Range("M136").Value = -12244.9995
$ws.Range("N136").Value = -15792.75

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 648.43475
$ws.Range("I3").Value = 686.0625
$ws.Range("J3").Value = 562.4286
$ws.Range("K3").Value = 686.0625
$ws.Range("L3").Value = 562.4286
$ws.Range("M3").Value = -572.0625
$ws.Range("N3").Value = -790.4286

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 100762.4
$ws.Range("I16").Value = 143587.28
$ws.Range("J16").Value = 837.6667
$ws.Range("K16").Value = 143587.28
$ws.Range("L16").Value = 837.6667
$ws.Range("M16").Value = -143300.28
$ws.Range("N16").Value = -1411.6667
# Row 31
$ws.Range("H31").Value = 4699.9893
$ws.Range("I31").Value = 844.4792
$ws.Range("J31").Value = 8723.130999999999
$ws.Range("K31").Value = 844.4792
$ws.Range("L31").Value = 8723.130999999999
$ws.Range("M31").Value = -549.4792
$ws.Range("N31").Value = -9313.130999999999
# Row 34
$ws.Range("H34").Value = 4699.9893
$ws.Range("I34").Value = 844.4792
$ws.Range("J34").Value = 8723.130999999999
$ws.Range("K34").Value = 844.4792
$ws.Range("L34").Value = 8723.130999999999
$ws.Range("M34").Value = -642.4792
$ws.Range("N34").Value = -9127.130999999999
# Row 58
$ws.Range("H58").Value = 10437.971
$ws.Range("I58").Value = 712.5238000000001
$ws.Range("J58").Value = 26148.309
$ws.Range("K58").Value = 712.5238000000001
$ws.Range("L58").Value = 26148.309
$ws.Range("M58").Value = -509.5238000000001
$ws.Range("N58").Value = -26554.309
# Row 113
$ws.Range("H113").Value = 100762.4
$ws.Range("I113").Value = 143587.28
$ws.Range("J113").Value = 837.6667
$ws.Range("K113").Value = 143587.28
$ws.Range("L113").Value = 837.6667
$ws.Range("M113").Value = -141417.28
$ws.Range("N113").Value = -5177.6667
# Row 122
$ws.Range("H122").Value = 23811316
$ws.Range("I122").Value = 47620196
$ws.Range("J122").Value = 2435
$ws.Range("K122").Value = 142860588
$ws.Range("L122").Value = 7305
$ws.Range("M122").Value = -142858138
$ws.Range("N122").Value = -12205
# Row 136
$ws.Range("H136").Value = 10437.971
$ws.Range("I136").Value = 712.5238000000001
$ws.Range("J136").Value = 26148.309
$ws.Range("K136").Value = 2137.5714
$ws.Range("L136").Value = 78444.927
$ws.Range("M136").Value = 412.4285999999997
$ws.Range("N136").Value = -83544.927

$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 3376.1428
$ws.Range("I98").Value = 5166.6665
$ws.Range("J98").Value = 2033.25
$ws.Range("K98").Value = 15499.9995
$ws.Range("L98").Value = 6099.75
$ws.Range("M98").Value = -14001.9995
$ws.Range("N98").Value = -9095.75
# Row 102
$ws.Range("H102").Value = 2999
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 2999
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 8997
$ws.Range("N102").Value = -13865
$ws.Range("M102").Value = $null
# Row 107
$ws.Range("H107").Value = 656.49396
$ws.Range("I107").Value = 325.30304
$ws.Range("J107").Value = 1942.2941
$ws.Range("K107").Value = 975.90912
$ws.Range("L107").Value = 5826.8823
$ws.Range("M107").Value = 944.09088
$ws.Range("N107").Value = -9666.882300000001

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 21117.31
$ws.Range("I70").Value = 27033.023
$ws.Range("J70").Value = 4158.933
$ws.Range("K70").Value = 27033.023
$ws.Range("L70").Value = 4158.933
$ws.Range("M70").Value = -26763.023
$ws.Range("N70").Value = -4698.933
# Row 73
$ws.Range("H73").Value = 21117.31
$ws.Range("I73").Value = 27033.023
$ws.Range("J73").Value = 4158.933
$ws.Range("K73").Value = 27033.023
$ws.Range("L73").Value = 4158.933
$ws.Range("M73").Value = -26097.023
$ws.Range("N73").Value = -6030.933
# Row 97
$ws.Range("H97").Value = 142858290
$ws.Range("I97").Value = 111112130
$ws.Range("J97").Value = 200001390
$ws.Range("K97").Value = 111112130
$ws.Range("L97").Value = 200001390
$ws.Range("M97").Value = -111111634
$ws.Range("N97").Value = -200002382
# Row 113
$ws.Range("H113").Value = 1027.7142
$ws.Range("I113").Value = 1039.6
$ws.Range("K113").Value = 1039.6
$ws.Range("M113").Value = 1130.4
# Row 132
$ws.Range("H132").Value = 43482936
$ws.Range("I132").Value = 58824880
$ws.Range("J132").Value = 14085.5
$ws.Range("K132").Value = 176474640
$ws.Range("L132").Value = 42256.5
$ws.Range("M132").Value = -176472110
$ws.Range("N132").Value = -47316.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 6588.8887
$ws.Range("I46").Value = 1150
$ws.Range("K46").Value = 1150
$ws.Range("M46").Value = -962
# Row 61
$ws.Range("H61").Value = 4502
$ws.Range("I61").Value = 4004
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 4004
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -3802
$ws.Range("N61").Value = -5404
# Row 113
$ws.Range("H113").Value = 4502
$ws.Range("I113").Value = 4004
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4004
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -1834
$ws.Range("N113").Value = -9340
# Row 133
$ws.Range("H133").Value = 20163
$ws.Range("J133").Value = 20163
$ws.Range("L133").Value = 20163
$ws.Range("N133").Value = -25223

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 11140.6
$ws.Range("I96").Value = 17534.334
$ws.Range("J96").Value = 1550
$ws.Range("K96").Value = 17534.334
$ws.Range("L96").Value = 1550
$ws.Range("M96").Value = -16161.334
$ws.Range("N96").Value = -4296
# Row 122
$ws.Range("H122").Value = 6043.5557
$ws.Range("I122").Value = 9598.091
$ws.Range("J122").Value = 457.85715
$ws.Range("K122").Value = 28794.273
$ws.Range("L122").Value = 1373.57145
$ws.Range("M122").Value = -26344.273
$ws.Range("N122").Value = -6273.571449999999
